$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed cells in existing rows (revised loc combination values) ---
$ws.Range("B1").Value = 0.443
$ws.Range("B2").Value = 0.885
$ws.Range("C2").Value = 0.003
$ws.Range("B3").Value = 0.852
$ws.Range("C3").Value = 0.003
$ws.Range("B4").Value = 0.439
$ws.Range("B5").Value = 0.442
$ws.Range("B6").Value = 1.652
$ws.Range("D6").Value = 0.6870000000000001
$ws.Range("E6").Value = 1
$ws.Range("B7").Value = 1.144
$ws.Range("C7").Value = 0.004
$ws.Range("D7").Value = 0.759
$ws.Range("E7").Value = 1
$ws.Range("B8").Value = 1.803
$ws.Range("C8").Value = 0.006
$ws.Range("D8").Value = 0.967
$ws.Range("B9").Value = 1.173
$ws.Range("C9").Value = 0.004
$ws.Range("D9").Value = 0.785
$ws.Range("E9").Value = 0.714
$ws.Range("B10").Value = 1.561
$ws.Range("C10").Value = 0.005
$ws.Range("D10").Value = 0.735
$ws.Range("B11").Value = 0.856
$ws.Range("C11").Value = 0.003
$ws.Range("D11").Value = 0.794
$ws.Range("B12").Value = 1.11
$ws.Range("C12").Value = 0.004
$ws.Range("B13").Value = 0.884
$ws.Range("B14").Value = 0.929
$ws.Range("B15").Value = 0.995
$ws.Range("D15").Value = 0.858
$ws.Range("E15").Value = 1
$ws.Range("B16").Value = 0.714
$ws.Range("B17").Value = 0.697
$ws.Range("B18").Value = 1.073
$ws.Range("C18").Value = 0.004
$ws.Range("D18").Value = 0.91
$ws.Range("B19").Value = 0.629
$ws.Range("B20").Value = 1.018
$ws.Range("D20").Value = 0.8139999999999999
$ws.Range("B21").Value = 0.511
$ws.Range("C21").Value = 0.002
$ws.Range("B22").Value = 1.249
$ws.Range("D22").Value = 0.726
$ws.Range("B23").Value = 1.168
$ws.Range("C23").Value = 0.004
$ws.Range("D23").Value = 0.847
$ws.Range("B24").Value = 1.484
$ws.Range("C24").Value = 0.005
$ws.Range("D24").Value = 0.895
$ws.Range("B25").Value = 1.367
$ws.Range("C25").Value = 0.005
$ws.Range("D25").Value = 0.908
$ws.Range("B26").Value = 1.116
$ws.Range("C26").Value = 0.004
$ws.Range("D26").Value = 0.547
$ws.Range("E26").Value = 0.571
$ws.Range("B27").Value = 1.365
$ws.Range("C27").Value = 0.005
$ws.Range("D27").Value = 0.905
$ws.Range("B28").Value = 1.466
$ws.Range("C28").Value = 0.005
$ws.Range("D28").Value = 0.901
$ws.Range("B29").Value = 0.852
$ws.Range("C29").Value = 0.003
$ws.Range("B30").Value = 1.667
$ws.Range("C30").Value = 0.006
$ws.Range("D30").Value = 0.902
$ws.Range("B31").Value = 0.6909999999999999
$ws.Range("D31").Value = 0.832
$ws.Range("B32").Value = 0.723
$ws.Range("D32").Value = 0.79
$ws.Range("B33").Value = 1.121
$ws.Range("C33").Value = 0.004
$ws.Range("D33").Value = 0.888
$ws.Range("E33").Value = 0.571
$ws.Range("B34").Value = 1.677
$ws.Range("C34").Value = 0.006
$ws.Range("D34").Value = 0.869
$ws.Range("B35").Value = 1.487
$ws.Range("C35").Value = 0.005
$ws.Range("D35").Value = 0.893
$ws.Range("E35").Value = 0.857
$ws.Range("A36").Value = "FY_4.png"
$ws.Range("B36").Value = 1.258
$ws.Range("C36").Value = 0.004
$ws.Range("D36").Value = 0.834
$ws.Range("E36").Value = 0.714
$ws.Range("F36").Value = "Fanny Yusuf"
$ws.Range("A37").Value = "TO_1.png"
$ws.Range("B37").Value = 0.777
$ws.Range("D37").Value = 0.806
$ws.Range("A38").Value = "TO_2.png"
$ws.Range("B38").Value = 1.053
$ws.Range("D38").Value = 0.857
$ws.Range("A39").Value = "TO_3.png"
$ws.Range("B39").Value = 0.8090000000000001
$ws.Range("C39").Value = 0.003
$ws.Range("D39").Value = 0.858
$ws.Range("E39").Value = 1
$ws.Range("F39").Value = "Tiara Oktavian"
$ws.Range("G39").Value = "Benar"
$ws.Range("A40").Value = "TO_4.png"
$ws.Range("B40").Value = 7.424
$ws.Range("C40").Value = 0.024
$ws.Range("D40").Value = 0.624
$ws.Range("E40").Value = 0.429
$ws.Range("F40").Value = "Tidak Diketahui"
$ws.Range("G40").Value = "Salah"
$ws.Range("A41").Value = "TO_5.png"
$ws.Range("B41").Value = 7.149
$ws.Range("C41").Value = 0.024
$ws.Range("D41").Value = 0.553
$ws.Range("E41").Value = 0.429
$ws.Range("F41").Value = "Tidak Diketahui"
$ws.Range("A42").Value = "TD_1.png"
$ws.Range("B42").Value = 3.795
$ws.Range("C42").Value = 0.013
$ws.Range("D42").Value = 0.456
$ws.Range("A43").Value = "TD_2.png"
$ws.Range("B43").Value = 5.005
$ws.Range("C43").Value = 0.017
$ws.Range("D43").Value = 0.488
$ws.Range("E43").Value = 0.571
$ws.Range("F43").Value = "Rafiqo Rapitasari"
$ws.Range("G43").Value = "Salah"
$ws.Range("A44").Value = "TD_3.png"
$ws.Range("B44").Value = 1.049
$ws.Range("D44").Value = 0.732
$ws.Range("E44").Value = 0.286
$ws.Range("F44").Value = "Tidak Diketahui"
$ws.Range("G44").Value = "Benar"

# --- Append new row 45 (TD_4.png) ---
$ws.Range("A45").Value = "TD_4.png"
$ws.Range("B45").Value = 0.994
$ws.Range("C45").Value = 0.003
$ws.Range("D45").Value = 0.718
$ws.Range("E45").Value = 0.286
$ws.Range("F45").Value = "Tidak Diketahui"
$ws.Range("G45").Value = "Benar"
